$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column C ("pembatalan_H03").
# Everything that was in C..AL shifts right to D..AM.
$ws.Columns("C").Insert()
$ws.Columns("C").ColumnWidth = 7.17

# Header for the newly inserted column
$ws.Range("C1").Value = "pembatalan_H03"

# Data fix-up for the two existing rows: "edit and delete data"
#  - row 2: status flag (col B) corrected from 3 -> 0, and the new
#    "pembatalan" flag (col C) is set to 1
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1

#  - row 3: status flag (col B) stays 3, new "pembatalan" flag (col C) set to 0
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 0

# Leave the selection where the user left it before saving
$ws.Range("K11").Select()
